$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (styles, row height) of the last existing row (13)
# down into the three new rows (14-16) so the new data matches the sheet's
# existing look (text columns formatted as Text, date column formatted as
# a short date, vertical-top alignment, row height, etc).
$ws.Range("A13:F13").Copy() | Out-Null
$ws.Range("A14:F16").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# Match row 13's explicit row height (the sheet default is 12.75; the
# data rows are set slightly taller at 13.05).
$ws.Rows.Item(14).RowHeight = $ws.Rows.Item(13).RowHeight
$ws.Rows.Item(15).RowHeight = $ws.Rows.Item(13).RowHeight
$ws.Rows.Item(16).RowHeight = $ws.Rows.Item(13).RowHeight

# Row 14: Paradise Smoothie Cafe
$ws.Range("A14").Value = "Paradise Smoothie Cafe "
$ws.Range("B14").Value = "Ballman, John W"
$ws.Range("C14").Value = "040"
$ws.Range("D14").Value = ""
$ws.Range("E14").Value = "0008329"

# Row 15: OWEN MEATS CORP
$ws.Range("A15").Value = "OWEN MEATS CORP"
$ws.Range("B15").Value = "Zigan, Gerald L"
$ws.Range("C15").Value = "030"
$ws.Range("D15").Value = ""
$ws.Range("E15").Value = "0008331"

# Row 16: HOLY FAMILY MARONITE CHURCH
$ws.Range("A16").Value = "HOLY FAMILY MARONITE CHURCH"
$ws.Range("B16").Value = "Bloch, Lea L"
$ws.Range("C16").Value = "003"
# D16 holds an actual date (unlike D14/D15), so borrow the format of an
# existing populated date cell (D7) instead of the blank-date format that
# was pasted down from D13.
$ws.Range("D7").Copy() | Out-Null
$ws.Range("D16").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0
# Assign the raw date serial (10/02/2025) rather than a .NET DateTime -
# writing a DateTime makes Excel stamp its own default short-date format
# over the cell, which would fork off a brand-new style instead of
# reusing the existing date style shared with D2/D3/D7.
$ws.Range("D16").Value = 45932
$ws.Range("E16").Value = "0004965"

$wb.Save()
